# Quarterly income-statement database update for Shavan (شاوان-پالایش نفت لاوان):
# drop the oldest reporting period column and append the newest one,
# shifting every period-dependent column one slot to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: financial-period column headers ---
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: matching "تاریخ انتشار" (publish date) labels ---
$ws.Range("D9").Value = "1400-11-05 (4)"
$ws.Range("E9").Value = "1401-04-18 (10)"
$ws.Range("F9").Value = "1401-04-21 (2)"
$ws.Range("G9").Value = "1401-08-30 (4)"
$ws.Range("H9").Value = "1401-10-29 (3)"
$ws.Range("I9").Value = "1402-02-27 (7)"
$ws.Range("J9").Value = "'1401-04-21"
$ws.Range("K9").Value = "1401-08-30 (2)"
$ws.Range("L9").Value = "'1401-10-29"
$ws.Range("M9").Value = "'1402-02-27"

# --- Rows 11-27: income-statement figures ---
# Row 11: فروش (Sales)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 75542234
$arr[0,1] = 139276849
$arr[0,2] = 65227719
$arr[0,3] = 143775382
$arr[0,4] = 226475563
$arr[0,5] = 326717314
$arr[0,6] = 132230271
$arr[0,7] = 269223523
$arr[0,8] = 372463034
$arr[0,9] = 454818108
$ws.Range("D11:M11").Value = $arr

# Row 12: بهای تمام شده کالای فروش رفته (COGS)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -68660675
$arr[0,1] = -115480712
$arr[0,2] = -55450974
$arr[0,3] = -121816846
$arr[0,4] = -198326745
$arr[0,5] = -294118552
$arr[0,6] = -109406194
$arr[0,7] = -230403284
$arr[0,8] = -331302090
$arr[0,9] = -415878679
$ws.Range("D12:M12").Value = $arr

# Row 13: سود (زیان) ناخالص (Gross profit)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 6881559
$arr[0,1] = 23796137
$arr[0,2] = 9776745
$arr[0,3] = 21958536
$arr[0,4] = 28148818
$arr[0,5] = 32598762
$arr[0,6] = 22824077
$arr[0,7] = 38820239
$arr[0,8] = 41160944
$arr[0,9] = 38939429
$ws.Range("D13:M13").Value = $arr

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -821207
$arr[0,1] = -1009498
$arr[0,2] = -400738
$arr[0,3] = -694083
$arr[0,4] = -959471
$arr[0,5] = -2192945
$arr[0,6] = -924523
$arr[0,7] = -1212005
$arr[0,8] = -1773624
$arr[0,9] = -3235786
$ws.Range("D14:M14").Value = $arr

# Row 15: هزینه کاهش ارزش دریافتنی ها (Impairment expense)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D15:M15").Value = $arr

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 3182849
$arr[0,1] = 1518153
$arr[0,2] = 970463
$arr[0,3] = 955644
$arr[0,4] = 842870
$arr[0,5] = 1618406
$arr[0,6] = 0
$arr[0,7] = 1329369
$arr[0,8] = 2279216
$arr[0,9] = 5729110
$ws.Range("D16:M16").Value = $arr

# Row 17: سود (زیان) عملیاتی (Operating profit)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 9243201
$arr[0,1] = 24304792
$arr[0,2] = 10346470
$arr[0,3] = 22220097
$arr[0,4] = 28032217
$arr[0,5] = 32024223
$arr[0,6] = 21899554
$arr[0,7] = 38937603
$arr[0,8] = 41666536
$arr[0,9] = 41432753
$ws.Range("D17:M17").Value = $arr

# Row 18: هزینه های مالی (Financial expenses)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -234281
$arr[0,1] = -484121
$arr[0,2] = -297126
$arr[0,3] = -517256
$arr[0,4] = -835536
$arr[0,5] = -1083194
$arr[0,6] = -266541
$arr[0,7] = -516557
$arr[0,8] = -878588
$arr[0,9] = -1292551
$ws.Range("D18:M18").Value = $arr

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 396059
$arr[0,1] = 71646
$arr[0,2] = 33240
$arr[0,3] = 204571
$arr[0,4] = 0
$arr[0,5] = 424271
$arr[0,6] = 40068
$arr[0,7] = 271664
$arr[0,8] = 537616
$arr[0,9] = 6306342
$ws.Range("D19:M19").Value = $arr

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 9404979
$arr[0,1] = 23892317
$arr[0,2] = 10082584
$arr[0,3] = 21907412
$arr[0,4] = 27196681
$arr[0,5] = 31365300
$arr[0,6] = 21673081
$arr[0,7] = 38692710
$arr[0,8] = 41325564
$arr[0,9] = 46446544
$ws.Range("D20:M20").Value = $arr

# Row 21: مالیات (Tax)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -1865219
$arr[0,1] = -3344763
$arr[0,2] = -2268582
$arr[0,3] = -4914014
$arr[0,4] = -6119253
$arr[0,5] = -2833005
$arr[0,6] = -5418270
$arr[0,7] = -8311008
$arr[0,8] = -6508776
$arr[0,9] = -4480019
$ws.Range("D21:M21").Value = $arr

# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit, continuing ops)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 7539760
$arr[0,1] = 20547554
$arr[0,2] = 7814002
$arr[0,3] = 16993398
$arr[0,4] = 21077428
$arr[0,5] = 28532295
$arr[0,6] = 16254811
$arr[0,7] = 30381702
$arr[0,8] = 34816788
$arr[0,9] = 41966525
$ws.Range("D22:M22").Value = $arr

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (Discontinued ops)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D23:M23").Value = $arr

# Row 24: سود (زیان) خالص (Net profit)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 7539760
$arr[0,1] = 20547554
$arr[0,2] = 7814002
$arr[0,3] = 16993398
$arr[0,4] = 21077428
$arr[0,5] = 28532295
$arr[0,6] = 16254811
$arr[0,7] = 30381702
$arr[0,8] = 34816788
$arr[0,9] = 41966525
$ws.Range("D24:M24").Value = $arr

# Row 25: سود هر سهم پس از کسر مالیات (EPS after tax)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 6594
$arr[0,1] = 2568
$arr[0,2] = 977
$arr[0,3] = 2124
$arr[0,4] = 2635
$arr[0,5] = 24953
$arr[0,6] = 2032
$arr[0,7] = 3798
$arr[0,8] = 4352
$arr[0,9] = 5246
$ws.Range("D25:M25").Value = $arr

# Row 26: سرمایه (Capital)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1143422
$arr[0,1] = 8000000
$arr[0,2] = 8000000
$arr[0,3] = 8000000
$arr[0,4] = 8000000
$arr[0,5] = 1143422
$arr[0,6] = 8000000
$arr[0,7] = 8000000
$arr[0,8] = 8000000
$arr[0,9] = 8000000
$ws.Range("D26:M26").Value = $arr

# Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS based on latest capital)
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 942
$arr[0,1] = 2568
$arr[0,2] = 977
$arr[0,3] = 2124
$arr[0,4] = 2635
$arr[0,5] = 3567
$arr[0,6] = 2032
$arr[0,7] = 3798
$arr[0,8] = 4352
$arr[0,9] = 5246
$ws.Range("D27:M27").Value = $arr

